$d = $word.ActiveDocument

# Replace all occurrences of "July 06, 2022" with "July 09, 2022"
$d.Content.Find.Execute("July 06, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "July 09, 2022", 2)

# Replace "September 04, 2022" with "September 07, 2022"
$d.Content.Find.Execute("September 04, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "September 07, 2022", 2)
